$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 46
$ws.Range("D46").Value = "X"
$ws.Range("E46").Value = "page.for"
$ws.Range("F46").Value = "VISION menu 4, listagem de telas existentes as quais foi adicionada uma nova tela chama olmsnp que chama uma subroutine OLMSNP"

# Row 47
$ws.Range("D47").Value = "X"
$ws.Range("E47").Value = "tstlog.for"
$ws.Range("F47").Value = "programa utilitario para testar escrita e leitura das transações no tmf (não chega a escrever no tmf mas sim usa tralog e logbuf em que passa o buff directemante)"

# Row 48
$ws.Range("D48").Value = "X"
$ws.Range("E48").Value = "kilsys.com"
$ws.Range("F48").Value = "script ou command procedure que serve para terminar abrutamente os processos relativos  ao sistema de Jogo Millennium "

# Row 49
$ws.Range("D49").Value = "DMQ_OLM_CONFIG.COM"
$ws.Range("E49").Value = " X"
$ws.Range("F49").Value = "script ou command procedure que serve para alterar o ficheiro de configuração do MessageQ que DMQ.INI e alterar qual é maquina que primaria e qual é a failover"

# Row 50
$ws.Range("D50").Value = "olmcommon.for"
$ws.Range("E50").Value = "X"
$ws.Range("F50").Value = "ficheiro de código fonte que contêm subroutines que servem para invocar RTL (Remove from The bottom of the List) e ABL (Add to the Bottom of the List) respectivamente."

# Apply style (wrap text, left/top alignment) to new F cells, matching the other
# cells in column F (column F already uses this style elsewhere on the sheet)
$ws.Range("F46:F50").HorizontalAlignment = -4131
$ws.Range("F46:F50").VerticalAlignment = -4160
$ws.Range("F46:F50").WrapText = $true

# Row heights (auto-calculated by Excel from the wrapped text content)
$ws.Rows("46:49").RowHeight = 30
$ws.Rows("50").RowHeight = 45

# Autofit columns D and E to the new (wider) content
$ws.Columns("D:E").AutoFit()

# Update selection / view (scroll so row 43 is at the top, select G47)
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G47").Select()
